$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset now covers two sending/target clusters ("ECs" and "M2") instead
# of just "M2", with "Cd84" as the ligand/receptor symbol throughout -- giving
# the full 2x2 cluster-pair matrix across rows 2-5.

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd84"
$ws.Range("C2").Value = "Cd84"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 35.70889233333333
$ws.Range("H2").Value = 107.126677
$ws.Range("I2").Value = 0.2383893603686217
$ws.Range("J2").Value = 0.2383893603686217
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 35.70889233333333
$ws.Range("N2").Value = 107.126677
$ws.Range("O2").Value = 0.2383893603686217
$ws.Range("P2").Value = 0.2383893603686217
$ws.Range("Q2").Value = 1275.124991673592
$ws.Range("R2").Value = 11476.12492506233
$ws.Range("S2").Value = 0.05682948713696057
$ws.Range("T2").Value = 0.0568294871369606

# Row 3: ECs -> M2
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd84"
$ws.Range("C3").Value = "Cd84"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 35.70889233333333
$ws.Range("H3").Value = 107.126677
$ws.Range("I3").Value = 0.2383893603686217
$ws.Range("J3").Value = 0.2383893603686217
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 114.083415
$ws.Range("N3").Value = 342.2502449999999
$ws.Range("O3").Value = 0.7616106396313783
$ws.Range("P3").Value = 0.7616106396313783
$ws.Range("Q3").Value = 4073.792383253985
$ws.Range("R3").Value = 36664.13144928586
$ws.Range("S3").Value = 0.1815598732316611
$ws.Range("T3").Value = 0.1815598732316611

# Row 4: M2 -> ECs
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Cd84"
$ws.Range("C4").Value = "Cd84"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 114.083415
$ws.Range("H4").Value = 342.2502449999999
$ws.Range("I4").Value = 0.7616106396313783
$ws.Range("J4").Value = 0.7616106396313783
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 35.70889233333333
$ws.Range("N4").Value = 107.126677
$ws.Range("O4").Value = 0.2383893603686217
$ws.Range("P4").Value = 0.2383893603686217
$ws.Range("Q4").Value = 4073.792383253985
$ws.Range("R4").Value = 36664.13144928586
$ws.Range("S4").Value = 0.1815598732316611
$ws.Range("T4").Value = 0.1815598732316611

# Row 5: M2 -> M2
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Cd84"
$ws.Range("C5").Value = "Cd84"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 114.083415
$ws.Range("H5").Value = 342.2502449999999
$ws.Range("I5").Value = 0.7616106396313783
$ws.Range("J5").Value = 0.7616106396313783
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 114.083415
$ws.Range("N5").Value = 342.2502449999999
$ws.Range("O5").Value = 0.7616106396313783
$ws.Range("P5").Value = 0.7616106396313783
$ws.Range("Q5").Value = 13015.02557806222
$ws.Range("R5").Value = 117135.23020256
$ws.Range("S5").Value = 0.5800507663997172
$ws.Range("T5").Value = 0.5800507663997172
